$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the three newly populated cells in the existing DNA007 (column E) rows
$ws.Range("E28").Value = 19
$ws.Range("E41").Value = 22
$ws.Range("E43").Value = 13

# Append the new Y-STR loci (Yfiler Plus) rows 50-54
$ws.Range("A50").Value = "DYS627"
$ws.Range("E50").Value = 21

$ws.Range("A51").Value = "DYS460"
$ws.Range("E51").Value = 11

$ws.Range("A52").Value = "DYS518"
$ws.Range("E52").Value = 37

$ws.Range("A53").Value = "DYS449"
$ws.Range("E53").Value = 30

$ws.Range("A54").Value = "DYF387S1"
$ws.Range("E54").Value = "35,37"

# Match the author's final view state (G44 selected)
$ws.Range("G44").Select()
